# Adds a "Time" column and recomputed "UPH" column (UserID / Qty / Time / UPH)
# to the PUTWALL PICKING and REGULAR PICK sheets, re-sorts the rows by the
# new UPH (descending), appends an "Average UPH" summary row, and highlights
# the header + summary rows with a light-blue fill.

$wb = $excel.ActiveWorkbook

$lightBlue = 15128749   # RGB(173,216,230) == #ADD8E6, packed as BGR for Excel COM

function Fill-Sheet {
    param(
        [string]$SheetName,
        [string]$QtyHeader,
        [object[,]]$Rows,
        [double]$AverageUph
    )

    # NOTE: positional args only -- this runtime does not bind named (-Foo)
    # parameters reliably, so Fill-Sheet is always invoked positionally below.
    $ws = $wb.Worksheets.Item($SheetName)

    $rowCount = $Rows.GetLength(0)
    $lastDataRow = 1 + $rowCount
    $avgRow = $lastDataRow + 1

    # ---- Header row ----
    $header = New-Object 'object[,]' 1,4
    $header[0,0] = "UserID"
    $header[0,1] = $QtyHeader
    $header[0,2] = "Time"
    $header[0,3] = "UPH"
    $ws.Range("A1:D1").Value = $header

    # ---- Data rows (already sorted by UPH desc) ----
    $data = New-Object 'object[,]' $rowCount,4
    for ($i = 0; $i -lt $rowCount; $i++) {
        $data[$i,0] = $Rows[$i,0]
        $data[$i,1] = $Rows[$i,1]
        $data[$i,2] = $Rows[$i,2]
        $data[$i,3] = $Rows[$i,3]
    }
    $ws.Range("A2:D$lastDataRow").Value = $data

    # ---- Average UPH summary row ----
    $ws.Range("A$avgRow").Value = "Average UPH"
    $ws.Range("B$avgRow").Value = ""
    $ws.Range("C$avgRow").Value = ""
    $ws.Range("D$avgRow").Value = $AverageUph

    # ---- Styling: light-blue fill on header row + summary row ----
    $ws.Range("A1:D1").Interior.Color = $lightBlue
    $ws.Range("A" + $avgRow + ":D" + $avgRow).Interior.Color = $lightBlue
}

# ---------------------------------------------------------------------------
# PUTWALL PICKING  (UserID, PutwallPickingQuantity, Time, UPH)
# ---------------------------------------------------------------------------
$putwallRows = New-Object 'object[,]' 20,4
$putwallData = @(
    @("DIAN4065.ENTRIALGO",    197, 44,  268.64),
    @("ANASTASIIA.MAKHTOUT",   276, 91,  181.98),
    @("ABHI4088.ABHISHEK",     286, 104, 165),
    @("KADE3054.ZONGO",        169, 67,  151.34),
    @("BOHD0676.KUSHLIAK",     173, 76,  136.58),
    @("LOWRHY-OTIENO.JAOKO",   120, 53,  135.85),
    @("TANI2739.HOSSAINISLA",  378, 167, 135.81),
    @("STAN9294.BAUER",        94,  46,  122.61),
    @("THIE6554.DIALLO",       351, 182, 115.71),
    @("WILDINE.JEUNE",         141, 88,  96.14),
    @("RAVI4279.THAKUR",       29,  10,  0),
    @("SEPIDEH.AZARIHASHJIN",  24,  6,   0),
    @("MDSAIFUL.ISLAM",        66,  19,  0),
    @("RARG046N.YEBOAH",       3,   0,   0),
    @("NESR2403.ATTALAH",      49,  18,  0),
    @("AGNE8120.CARUTH",       12,  19,  0),
    @("LOANA.MBONGO",          1,   0,   0),
    @("HARJ4282.SINGH",        51,  17,  0),
    @("DEVI789.SINGH",         46,  15,  0),
    @("YATI0689.YATIN",        39,  9,   0)
)
for ($i = 0; $i -lt $putwallData.Count; $i++) {
    $putwallRows[$i,0] = $putwallData[$i][0]
    $putwallRows[$i,1] = $putwallData[$i][1]
    $putwallRows[$i,2] = $putwallData[$i][2]
    $putwallRows[$i,3] = $putwallData[$i][3]
}
Fill-Sheet "PUTWALL PICKING" "PutwallPickingQuantity" $putwallRows 150.97

# ---------------------------------------------------------------------------
# REGULAR PICK  (UserID, RegularPickQuantity, Time, UPH)
# ---------------------------------------------------------------------------
$regularRows = New-Object 'object[,]' 24,4
$regularData = @(
    @("DIAN4065.ENTRIALGO",     179, 38,  282.63),
    @("BOHD0676.KUSHLIAK",      148, 110, 80.73),
    @("WILDINE.JEUNE",          22,  43,  30.7),
    @("SEPIDEH.AZARIHASHJIN",   14,  34,  24.71),
    @("AGNE8120.CARUTH",        8,   3,   0),
    @("NESR2403.ATTALAH",       5,   0,   0),
    @("ZAHIDGUL.MINHAS",        17,  19,  0),
    @("THIE6554.DIALLO",        17,  22,  0),
    @("TANI2739.HOSSAINISLA",   2,   0,   0),
    @("SURESH.DHAWAN",          2,   0,   0),
    @("STAN9294.BAUER",         18,  20,  0),
    @("RAVI4279.THAKUR",        2,   0,   0),
    @("RARG046N.YEBOAH",        20,  1,   0),
    @("MARI882N.ABDELKADER",    17,  15,  0),
    @("ANASTASIIA.MAKHTOUT",    17,  4,   0),
    @("LOWRHY-OTIENO.JAOKO",    35,  23,  0),
    @("LOANA.MBONGO",           1,   0,   0),
    @("KHINEHAYMAR.THAUNG",     1,   0,   0),
    @("KADE3054.ZONGO",         1,   0,   0),
    @("JEEW9554.SITUMUDALIG",   5,   3,   0),
    @("HARJ4282.SINGH",         2,   0,   0),
    @("DEVI789.SINGH",          1,   0,   0),
    @("ARJUNBHAI.PATEL",        27,  7,   0),
    @("ZAKI0190.PHILLIPHORS",   9,   12,  0)
)
for ($i = 0; $i -lt $regularData.Count; $i++) {
    $regularRows[$i,0] = $regularData[$i][0]
    $regularRows[$i,1] = $regularData[$i][1]
    $regularRows[$i,2] = $regularData[$i][2]
    $regularRows[$i,3] = $regularData[$i][3]
}
Fill-Sheet "REGULAR PICK" "RegularPickQuantity" $regularRows 104.69

Write-Output "edit complete"
